# Update Icos-Icosl.xlsx with newly computed TPM-based values.
# Sending-cluster set shrinks from {ECs, FAPs, MuSCs} to {FAPs, MuSCs}
# (ECs remains only as a Target cluster), so rows 2-7 get new content
# and rows 8-10 (old MuSCs-as-sender block) are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) is untouched by this update - only the data block changes.

$rows = @(
    @("FAPs","Icos","Icosl","ECs",1,0.3333333333333333,0.107019,0.321057,0.2365534751714524,0.2365534751714524,3,1,0.5407596666666666,1.622279,0.03618231591230665,0.03618231591230665,0.05787155876699999,0.5208440289029999,0.00855905256880748,0.00855905256880748),
    @("FAPs","Icos","Icosl","FAPs",1,0.3333333333333333,0.107019,0.321057,0.2365534751714524,0.2365534751714524,3,1,11.32416633333333,33.972499,0.7577017832003755,0.7577017832003754,1.211900956827,10.907108611443,0.1792369899596553,0.1792369899596552),
    @("FAPs","Icos","Icosl","MuSCs",1,0.3333333333333333,0.107019,0.321057,0.2365534751714524,0.2365534751714524,3,1,3.080487333333333,9.241461999999999,0.2061159008873179,0.2061159008873179,0.3296706739259999,2.967036065333999,0.04875743264298971,0.04875743264298971),
    @("MuSCs","Icos","Icosl","ECs",3,1,0.3453903333333333,1.036171,0.7634465248285476,0.7634465248285476,3,1,0.5407596666666666,1.622279,0.03618231591230665,0.03618231591230665,0.1867731615232222,1.680958453709,0.02762326334349917,0.02762326334349917),
    @("MuSCs","Icos","Icosl","FAPs",3,1,0.3453903333333333,1.036171,0.7634465248285476,0.7634465248285476,3,1,11.32416633333333,33.972499,0.7577017832003755,0.7577017832003754,3.91125758459211,35.201318261329,0.5784647932407203,0.5784647932407202),
    @("MuSCs","Icos","Icosl","MuSCs",3,1,0.3453903333333333,1.036171,0.7634465248285476,0.7634465248285476,3,1,3.080487333333333,9.241461999999999,0.2061159008873179,0.2061159008873179,1.063970546889111,9.575734922001997,0.1573584682443282,0.1573584682443282)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowVals = $rows[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $rowVals[$c]
    }
}

# Drop the now-unused old rows 8-10 (used to carry the MuSCs-as-sender /
# (ECs,FAPs,MuSCs)-as-target block) so the sheet shrinks to A1:T7.
$ws.Range("A8:T10").Delete()
